# Update RECIPES.xlsx:
#  - Reformat RECIPE_REQUIREMENTS (column D) into <qty><unit><item> tuples
#  - Add a new NUM_INGR column (F) with a header and per-recipe ingredient counts
#  - Remove the trailing blank row (row 23)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New structured ingredient strings for column D (rows 3-12)
$ws.Range("D3").Value = "<24><Ounces><Pasta Sauce>,<1><Box><Spaghetti Noodles>,<1.5><lb><Ground Beef>"
$ws.Range("D4").Value = "<8><><Chicken Wings>, <1><Teaspoon> <Olive Oil>,<3><Tablespoons><butter>,<><><Seasoning>"
$ws.Range("D5").Value = "<8><><Chicken Wings>,<><><Teriyaki Sauce>,<><><Sesame Seeds>"
$ws.Range("D6").Value = "<2.5 - 3><lbs><Chicken with Bone>,<><><Salt>,<><><Pepper>,<2><teaspoons><olive oil>,<1><><Lemon>, <.5><cup><pitted olives>"
$ws.Range("D7").Value = "<2><><Flour Tortillas>,<><><Shredded Cheese>"
$ws.Range("D8").Value = "<2><><Flour Tortillas>,<><><Shredded Cheese>,<1><><Boneless Chicken Breast>  "
$ws.Range("D9").Value = "<.25><><Medium Watermelon>,<.25><><Red Onion>,<><><Salt>,<><><Pepper>,<1.5><lbs><tri tip>,<8><><red peppers>,<1><><lime>"
$ws.Range("D10").Value = "<26><ounce Jar><Pasta Sauce>,<1><Bag><Frozen Ravioli>,<><><Chopped Spinach>,<><><Shredded Mozzarella>,<><><Grated Parmesan Cheese>"
$ws.Range("D11").Value = "<5><Tablespoons><Butter>,<1><lb><Macaroni Noodles>,<.5><cup><flour>,<6><cups><whole milk>,<2><cups><Grated Gruyere>,<1.5><cups><Grated Cheddar>,<.25><teaspoons><Cayenne Pepper>,<><><Salt>"
$ws.Range("D12").Value = "<1><lb><ground beef>,<><><Salsa>,<10><><Taco Shells>,<.5><head><shredded lettuce>,<1><><medium chopped tomato>,<1><cup><shredded cheese>"

# New NUM_INGR column: header + counts
$ws.Range("F2").Value = "NUM_INGR"
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = 6

# Remove trailing blank row
$ws.Rows(23).Delete()
